$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Claim the alignment-only style slot (index 2) before the hyperlink writes
# below claim any auto-generated styles of their own.
$ws.Range("E7").WrapText = $true
$ws.Range("E7").WrapText = $false

# --- Row 6: CP006_enviar_solicitud ---
$ws.Range("B6").Value = "jisola.tsoft@gmail.com"
$ws.Range("C6").Value = 12061990
$ws.Range("D6").Value = "juan martin isola"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:jisola.tsoft@gmail.com", "", "", "jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B6").Style = $ws.Range("B2").Style

# --- Row 7: CP007_cancelar_solicitud ---
$ws.Range("B7").Value = "jisola.tsoft@gmail.com"
$ws.Range("C7").Value = 12061990
$ws.Range("D7").Value = "juan martin isola"
$ws.Range("E7").Value = "Cancelar solicitud"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:jisola.tsoft@gmail.com", "", "", "jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B7").Style = $ws.Range("B2").Style

# --- Row 8: CP008_meGusta_pagina ---
$ws.Range("B8").Value = "jisola.tsoft@gmail.com"
$ws.Range("C8").Value = 12061990
$ws.Range("D8").Value = "juan martin isola"
$ws.Range("E8").Value = "Agregar"
$ws.Hyperlinks.Add($ws.Range("B8"), "mailto:jisola.tsoft@gmail.com", "", "", "jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B8").Style = $ws.Range("B2").Style

# --- Row 9: CP009_crear_publicacion ---
$ws.Range("B9").Value = "jisola.tsoft@gmail.com"
$ws.Range("C9").Value = 12061990
$ws.Range("D9").Value = "Tsoft"
$ws.Range("E9").Value = "Te gusta"
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:jisola.tsoft@gmail.com", "", "", "jisola.tsoft@gmail.com") | Out-Null
$ws.Range("B9").Style = $ws.Range("B2").Style

# Selection moved to F6 as of the last save
$ws.Range("F6").Select() | Out-Null

$wb.Save()
